$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos list values
$ws.Range("D2").Value = "37.467.36"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.021.67"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'253.21"
$ws.Range("E5").Value = "  +2.75%  "
$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'56.67"
$ws.Range("E8").Value = "  -8.13%  "
$ws.Range("D9").Value = "'0.381"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("E10").Value = "  -3.56%  "
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("D12").Value = "'14.50"
$ws.Range("E12").Value = "  -3.08%  "
$ws.Range("D13").Value = "2.321.41"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").Value = "'0.811"
$ws.Range("E14").Value = "  -4.94%  "
$ws.Range("D15").Value = "'20.96"
$ws.Range("E15").Value = "  -8.64%  "
$ws.Range("E16").Value = "  -2.54%  "
$ws.Range("D17").Value = "2.063.87"
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("D18").Value = "37.348.74"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "'69.47"
$ws.Range("E19").Value = "  -1.68%  "
$ws.Range("E20").Value = "  -2.98%  "
$ws.Range("D21").Value = "'5.18"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "'227.71"
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("E24").Value = "  +2.92%  "
$ws.Range("E25").Value = "  -2.23%  "
$ws.Range("D26").Value = "'163.12"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "'9.00"
$ws.Range("E27").Value = "  -4.36%  "
$ws.Range("D28").Value = "'19.74"
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("D29").Value = "'0.130"
$ws.Range("E29").Value = "  -9.61%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("D32").Value = "'0.0670"
$ws.Range("E32").Value = "  +6.73%  "
$ws.Range("E33").Value = "  -4.08%  "
$ws.Range("D34").Value = "'4.52"
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("E35").Value = "  +2.38%  "
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'3.41"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.82"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").Value = "'5.29"
$ws.Range("E39").Value = "  -5.05%  "
$ws.Range("E40").Value = "  +3.06%  "
$ws.Range("D41").Value = "'0.0961"
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("D43").Value = "'0.0214"
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D44").Value = "1.408.54"
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("D45").Value = "'15.87"
$ws.Range("E45").Value = "  -5.28%  "
$ws.Range("D46").Value = "'90.30"
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("D48").Value = "'7.27"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("D49").Value = "'2.86"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "2.213.26"
$ws.Range("E51").Value = "  +0.59%  "
